# Adds a new vintage column CF (date 45986) to the real-time GDP matrix,
# mirroring the existing CE column pattern, and appends a new trailing
# row (141) for the next vintage date (45976), matching the staircase
# structure of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column CF: header cell (row 1), styled like CE1 (date format) ---
$ws.Range("CE1").Copy($ws.Range("CF1"))
$ws.Range("CF1").Value = 45986

# --- Rows 3 through 138: CF gets the same value as CE in that row ---
for ($r = 3; $r -le 138; $r++) {
    $ws.Range("CE$r").Copy($ws.Range("CF$r"))
}

# --- Row 139: CF has its own (revised) value, not a straight copy of CE ---
$ws.Range("CF139").Value = -0.2099036351493167

# --- Row 140: CF is a fresh 0 value (first estimate for that vintage) ---
$ws.Range("CF140").Value = 0

# --- Row 141: brand-new trailing row for the next vintage date, styled like A140 ---
$ws.Range("A140").Copy($ws.Range("A141"))
$ws.Range("A141").Value = 45976
